# Redefine definition of financial cost as effective cost/loan ratio
#
# The visible "oc_reg" sheet pulls every number through formulas such as
# =[1]oc_reg!C5 that reference a cached external link (xl/externalLinks/
# externalLink1.xml, itself fed from reg_results/oc_reg.csv, on disk next
# to the workbook). The author re-ran the underlying regression (changing
# how "financing cost" is computed), which changed the csv, refreshed the
# external link cache, and rippled into the handful of display cells below.
#
# That source csv isn't reachable from this sandbox (it's an external,
# disk-based link - there's no Excel object for editing a link's cached
# values directly), so we reproduce the refreshed numbers directly on the
# cells that actually changed, leaving every other (unchanged) linked
# formula cell untouched.
#
# Every one of these cells holds text (the external csv's numbers are
# cached as strings, e.g. "-125.0", "(0.072)"), so each literal is entered
# with a leading apostrophe to force text entry - exactly what typing the
# same digits into Excel's UI would do - instead of letting AutoDetect
# coerce it to a number (which would silently drop the formatting, e.g.
# "-125.0" -> -125, "(0.072)" -> -0.072).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $text) {
    $ws.Range($addr).Formula = "'" + $text
}

# Row 5 (OC (dummy) coefficient row)
Set-Text "B5" "-0.16"
Set-Text "D5" "0.039"
Set-Text "E5" "0.010"
Set-Text "F5" "-129.5"
Set-Text "G5" "-125.0"

# Row 6 (std. error row under OC (dummy))
Set-Text "B6" "(0.072)"
Set-Text "F6" "(33.3)"
Set-Text "G6" "(34.6)"

# Row 7 (Constant coefficient row)
Set-Text "B7" "0.27"
Set-Text "C7" "0.24"
Set-Text "D7" "0.61"
Set-Text "E7" "0.49"
Set-Text "F7" "-307.2"
Set-Text "G7" "-289.3"

# Row 8 (std. error row under Constant)
Set-Text "B8" "(0.086)"
Set-Text "E8" "(0.18)"
Set-Text "F8" "(62.8)"
Set-Text "G8" "(64.3)"

# Row 10 (Observations)
Set-Text "B10" "1040"
Set-Text "D10" "845"
Set-Text "F10" "1265"
Set-Text "G10" "1265"

# Row 11 (R-sq)
Set-Text "B11" "0.137"
Set-Text "D11" "0.189"
Set-Text "F11" "0.042"
Set-Text "G11" "0.043"

# Row 12 (Dep. Var. Mean)
Set-Text "F12" "-323.2"
